# Updated cryptos list on Fri Jul  7 20:41:40 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows. The Price column stores values as literal text (e.g.
# "30.257.28", "1.000") rather than numbers, so each Price update is
# entered with a leading apostrophe to force text entry instead of letting
# Excel auto-convert it to a number; the cell's Style is then reset to
# "Normal" so no stray "text format" styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.252.82"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.863.98"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'234.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = "'0.4674"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = "'0.2835"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = "'0.06510"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').Value = "'21.26"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07851"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = "'97.07"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').Value = "'1.874.20"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = "'5.084"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = "'0.6719"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = "'279.72"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').Value = "'30.256.69"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = "'1.001"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = "'5.474"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').Value = "'12.65"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = "'2.117.94"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = "'0.000007257"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = "'6.140"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = "'9.177"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.25%  '
$ws.Range('D26').Value = "'164.81"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('D27').Value = "'19.08"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = "'1.923"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.58%  '
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = "'0.09636"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D31').Value = "'4.383"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').Value = "'1.475"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').Value = "'4.089"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').Value = "'0.04698"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').Value = "'1.117"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.74%  '
$ws.Range('D36').Value = "'0.7043"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = "'2.727"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').Value = "'0.01848"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('D39').Value = "'2.534"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = "'6.244"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.21%  '
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = "'1.939"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').Value = "'0.8439"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').Value = "'0.4164"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').Value = "'0.9998"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = "'103.66"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').Value = "'7.165"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.54%  '
$ws.Range('D48').Value = "'9.223"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').Value = "'936.17"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.11%  '
$ws.Range('D50').Value = "'34.01"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E51').Value = '  -1.87%  '
